# Insert one new price-record row for "Macroferia Regional de Talca - Cilantro"
# (weekly fruit/vegetable price refresh). This shifts the existing rows
# 99-120 down to 100-121 and fills the freed-up row 99 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 99..120 down to 100..121, opening up a blank row 99.
$ws.Rows.Item(99).Insert()

# Populate the new row 99 with the new weekly record.
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 45173
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112040
$ws.Range("G99").Value = "Cilantro"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 300
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = 8000
$ws.Range("N99").Value = "$/caja 36 atados"
$ws.Range("O99").Value = "Región Metropolitana"
$ws.Range("P99").Value = 222
$ws.Range("Q99").Value = 36
$ws.Range("R99").Value = "Hortaliza"
